# Apply the "New country level SLAND values GCB2023" update.
# The header row is restructured (SLAND_corrected/SLAND_dor removed, several
# columns re-ordered, and new *_percent columns appended), and all data rows
# (BRICS, LDC, OECD, Other) get refreshed values, extending the used range
# from A1:Q5 to A1:R5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 'ELUC'
$ws.Range("F1").Value = 'EFOS'
$ws.Range("G1").Value = 'F_ab'
$ws.Range("H1").Value = 'F_ac'
$ws.Range("I1").Value = 'F_abc'
$ws.Range("J1").Value = 'continent'
$ws.Range("K1").Value = 'climate_zone'
$ws.Range("L1").Value = 'CSCC_percent'
$ws.Range("M1").Value = 'SLAND_percent'
$ws.Range("N1").Value = 'ELUC_percent'
$ws.Range("O1").Value = 'EFOS_percent'
$ws.Range("P1").Value = 'F_ab_percent'
$ws.Range("Q1").Value = 'F_ac_percent'
$ws.Range("R1").Value = 'F_abc_percent'
$ws.Range("C2").Value = 296.7258181216166
$ws.Range("D2").Value = 1.126923322677612
$ws.Range("E2").Value = -0.3258713675000001
$ws.Range("F2").Value = -4.536557005036167
$ws.Range("G2").Value = 0.8010519327531108
$ws.Range("H2").Value = -4.862428372536167
$ws.Range("I2").Value = -3.735505072283056
$ws.Range("J2").Value = 'AsiaLatin America and the CaribbeanAsiaAfricaAfricaAsiaAsiaAsiaAfrica'
$ws.Range("K2").Value = '214211252'
$ws.Range("L2").Value = 44.79287594714138
$ws.Range("M2").Value = 33.44101715087891
$ws.Range("N2").Value = 27.74531822575011
$ws.Range("O2").Value = 48.39513225354158
$ws.Range("P2").Value = 36.48817090522742
$ws.Range("Q2").Value = 46.09590566623369
$ws.Range("R2").Value = 52.03653362319387
$ws.Range("C3").Value = 67.87593501224863
$ws.Range("D3").Value = 0.5636439919471741
$ws.Range("E3").Value = -0.33841027275
$ws.Range("F3").Value = -0.04997149098791279
$ws.Range("G3").Value = 0.2252336994837912
$ws.Range("H3").Value = -0.3883817637379128
$ws.Range("I3").Value = 0.1752622084958784
$ws.Range("J3").Value = 'AsiaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAfricaLatin America and the CaribbeanAsiaAsiaAfricaAfricaAfricaAfricaAsiaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAsiaAfricaAfricaAfricaAfrica'
$ws.Range("K3").Value = '4111111122111111121211121112112111112'
$ws.Range("L3").Value = 10.24635589867598
$ws.Range("M3").Value = 16.72591781616211
$ws.Range("N3").Value = 28.81290485980374
$ws.Range("O3").Value = 0.5330864161041047
$ws.Range("P3").Value = 10.25946681401262
$ws.Range("Q3").Value = 3.681865885133947
$ws.Range("R3").Value = -2.441447040974577
$ws.Range("C4").Value = 40.27682264769349
$ws.Range("D4").Value = 0.8203729391098022
$ws.Range("E4").Value = -0.02166024025
$ws.Range("F4").Value = -3.388613410369376
$ws.Range("G4").Value = 0.7987126569866662
$ws.Range("H4").Value = -3.410273650619376
$ws.Range("I4").Value = -2.58990075338271
$ws.Range("J4").Value = 'OceaniaEuropeEuropeNorth AmericaEuropeLatin America and the CaribbeanLatin America and the CaribbeanLatin America and the CaribbeanEuropeEuropeEuropeEuropeEuropeEuropeEuropeEuropeEuropeEuropeEuropeEuropeAsiaEuropeAsiaAsiaEuropeEuropeEuropeLatin America and the CaribbeanEuropeEuropeOceaniaEuropeEuropeEuropeEuropeEuropeAsiaNorth America'
$ws.Range("K4").Value = '23353211333344333335233343413534333443'
$ws.Range("L4").Value = 6.080073287265186
$ws.Range("M4").Value = 24.34425163269043
$ws.Range("N4").Value = 1.844194729941872
$ws.Range("O4").Value = 36.14908706512399
$ws.Range("P4").Value = 36.38161614832529
$ws.Range("Q4").Value = 32.32945360858864
$ws.Range("R4").Value = 36.07797473870554
$ws.Range("B5").Value = 'ALBARGARMAZEBGDBGRBHSBIHBLRBLZBOLBRNBTNBWACIVCMRCOGCOKCPVCUBCYPDOMDZAECUFJIGABGEOGHAGNQGTMGUYHNDHRVIDNIRQJAMJORKAZKENKGZKIRKWTLBNLBYLKAMARMDAMKDMNEMNGMRTMUSMYSNAMNGANICNPLOMNPAKPANPERPHLPNGPRKPRYPSEQATROUSAUSDNSLBSLVSOMSOMSRBSURSYRTHATJKTKMTTOTUNUKRURYUZBVCTVENVNMVUTWSMYEM'
$ws.Range("C5").Value = 257.5611802005601
$ws.Range("D5").Value = 0.8589435219764709
$ws.Range("E5").Value = -0.4885675082500001
$ws.Range("F5").Value = -1.398852530214742
$ws.Range("G5").Value = 0.3703760276577128
$ws.Range("H5").Value = -1.887420038464742
$ws.Range("I5").Value = -1.028476502557029
$ws.Range("J5").Value = 'EuropeLatin America and the CaribbeanAsiaAsiaAsiaEuropeLatin America and the CaribbeanEuropeEuropeLatin America and the CaribbeanLatin America and the CaribbeanAsiaAsiaAfricaAfricaAfricaAfricaOceaniaAfricaLatin America and the CaribbeanAsiaLatin America and the CaribbeanAfricaLatin America and the CaribbeanOceaniaAfricaAsiaAfricaAfricaLatin America and the CaribbeanLatin America and the CaribbeanLatin America and the CaribbeanEuropeAsiaAsiaLatin America and the CaribbeanAsiaAsiaAfricaAsiaOceaniaAsiaAsiaAfricaAsiaAfricaEuropeEuropeEuropeAsiaAfricaAfricaAsiaAfricaAfricaLatin America and the CaribbeanAsiaAsiaAsiaLatin America and the CaribbeanLatin America and the CaribbeanAsiaOceaniaAsiaLatin America and the CaribbeanAsiaAsiaEuropeAsiaAfricaOceaniaLatin America and the CaribbeanAfricaAfricaEuropeLatin America and the CaribbeanAsiaAsiaAsiaAsiaLatin America and the CaribbeanAfricaEuropeLatin America and the CaribbeanAsiaLatin America and the CaribbeanLatin America and the CaribbeanAsiaOceaniaOceaniaAsia'
$ws.Range("K5").Value = '3144131341114211112131211141111131212414123212433421121142211114122322112231214212414111112'
$ws.Range("L5").Value = 38.88069486691744
$ws.Range("M5").Value = 25.48882102966309
$ws.Range("N5").Value = 41.59758218450427
$ws.Range("O5").Value = 14.92269426523033
$ws.Range("P5").Value = 16.87074613243467
$ws.Range("Q5").Value = 17.89277484004373
$ws.Range("R5").Value = 14.32693867907519
